# Slide 3 ("Что нужно, чтобы использовать сервис"), the bulleted placeholder
# shape ("Объект 2"): merge the first bullet's two runs into one run, and
# insert the word "к" into the second bullet (between "доступом" and "Wi-Fi"),
# splitting "Телефон с доступом " into "Телефон с " + "доступом к Wi".
#
# Note: TextRange/Paragraphs objects returned by this object model are live
# (their Start/Length track the text as it is edited), so each step below
# re-reads the paragraph right before using it instead of relying on sizes
# computed earlier.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item("Объект 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Bullet 1: "Микроконтроллер " + "ESP32" -> single run "Микроконтроллер ESP32"
$oldLead = "Микроконтроллер "
$newRun1 = "Микроконтроллер ESP32"

$para1 = $tr.Paragraphs(1, 1)

# Clear the leading run's text so the (now empty) run is dropped and the
# "ESP32" run shifts up to the start of the paragraph.
$lead = $tr.Characters($para1.Start, $oldLead.Length)
$lead.Text = ""

# Re-fetch the paragraph (now just "ESP32" + the paragraph mark) and replace
# it wholesale with the merged text so it keeps that remaining run's single
# set of character formatting.
$para1 = $tr.Paragraphs(1, 1)
$tailLen = $para1.Length - 1
$merged = $tr.Characters($para1.Start, $tailLen)
$merged.Text = $newRun1

# --- Bullet 2: "Телефон с доступом " + "Wi" + "-Fi"
#           -> "Телефон с " + "доступом к Wi" + "-Fi"
$oldFirstRun  = "Телефон с доступом "
$newFirstRun  = "Телефон с "
$oldSecondRun = "Wi"
$newSecondRun = "доступом к Wi"

$para2 = $tr.Paragraphs(2, 1)

# Shrink the first run down to "Телефон с ".
$firstRun = $tr.Characters($para2.Start, $oldFirstRun.Length)
$firstRun.Text = $newFirstRun

# The next run ("Wi") now starts right after "Телефон с "; replace it
# wholesale with "доступом к Wi".
$para2 = $tr.Paragraphs(2, 1)
$secondRunStart = $para2.Start + $newFirstRun.Length
$secondRun = $tr.Characters($secondRunStart, $oldSecondRun.Length)
$secondRun.Text = $newSecondRun
